$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $text) {
    $c = $sheet.Range($addr)
    $escaped = $text -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue $ws "D2" "288.31"
Set-TextValue $ws "E2" "-1.07%"
Set-TextValue $ws "D3" "31.04"
Set-TextValue $ws "E3" "1.31%"
Set-TextValue $ws "D4" "4.922"
Set-TextValue $ws "E4" "-0.56%"
Set-TextValue $ws "D5" "0.07339"
Set-TextValue $ws "E5" "1.74%"
Set-TextValue $ws "D6" "2.199"
Set-TextValue $ws "E6" "19.23%"
Set-TextValue $ws "D7" "7.718"
Set-TextValue $ws "E7" "0.45%"
Set-TextValue $ws "D8" "3.732"
Set-TextValue $ws "E8" "-0.80%"
Set-TextValue $ws "D9" "0.9021"
Set-TextValue $ws "E9" "0.49%"
Set-TextValue $ws "D10" "0.09158"
Set-TextValue $ws "E10" "18.80%"
Set-TextValue $ws "E11" "1.02%"
Set-TextValue $ws "D12" "0.08213"
Set-TextValue $ws "E12" "1.82%"
Set-TextValue $ws "D13" "0.03121"
Set-TextValue $ws "E13" "2.78%"
Set-TextValue $ws "D14" "0.09948"
Set-TextValue $ws "E14" "-0.70%"
Set-TextValue $ws "D15" "0.001499"
Set-TextValue $ws "E15" "0.20%"
Set-TextValue $ws "D16" "0.005738"
Set-TextValue $ws "E16" "0.20%"
Set-TextValue $ws "D17" "3.524"
Set-TextValue $ws "E17" "1.56%"
Set-TextValue $ws "D18" "2.066"
Set-TextValue $ws "E18" "-0.81%"
Set-TextValue $ws "E20" "0.35%"
Set-TextValue $ws "D21" "4.206"
Set-TextValue $ws "E21" "3.94%"
Set-TextValue $ws "D23" "0.04537"
Set-TextValue $ws "E23" "0.76%"
Set-TextValue $ws "D24" "0.001210"
Set-TextValue $ws "E24" "-0.55%"
Set-TextValue $ws "D25" "0.004155"
Set-TextValue $ws "E25" "3.56%"
Set-TextValue $ws "E26" "4.04%"
Set-TextValue $ws "D27" "0.0003397"
Set-TextValue $ws "E27" "-95.47%"
Set-TextValue $ws "D39" "0.01567"
Set-TextValue $ws "E39" "-2.55%"
Set-TextValue $ws "D40" "0.04439"
Set-TextValue $ws "E40" "0.50%"
Set-TextValue $ws "D41" "0.007276"
Set-TextValue $ws "E41" "-0.05%"
Set-TextValue $ws "D42" "0.008979"
Set-TextValue $ws "E42" "-9.57%"
Set-TextValue $ws "D43" "0.1327"
Set-TextValue $ws "E43" "1.49%"
Set-TextValue $ws "D44" "0.002233"
Set-TextValue $ws "E44" "11.18%"
Set-TextValue $ws "D45" "0.009096"
Set-TextValue $ws "E45" "-4.40%"
Set-TextValue $ws "D46" "0.00006123"
Set-TextValue $ws "E46" "2.67%"
Set-TextValue $ws "D48" "2.298"
Set-TextValue $ws "E48" "2.32%"
Set-TextValue $ws "D49" "0.002001"
Set-TextValue $ws "E49" "-33.33%"
Set-TextValue $ws "D50" "0.00002101"
Set-TextValue $ws "D51" "0.0002001"
